$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between C10 (Joueur 2 / Distribution tour 4) and C12 (Joueur 2 / Distribution tour 6)
$ws.Range("C10").Value = "Somnifères"
$ws.Range("C12").Value = "Gourde"

# Change E13 (Joueur 4 / Distribution tour 7) to a new item
$ws.Range("E13").Value = "Brosse à WC"

# Update the active cell selection to match the saved view state
$ws.Range("H19").Select()
